# DC70_FAV_SUAP.docx edits:
#  1. "Allo" -> "Al" (salutation to SUAP/SUE table cell)
#  2. "Sportello unico per le attività produttive" -> "SUAP/SUE"
#  3. " SUAP " -> " SUAP/SUE " (In risposta a nota ...)
#  4. "ratica SUAP n°" -> "ratica SUAP/SUE n°"
#  5. Remove the whole "Si premette che il SUAP riceve ..." paragraph
#  6. Footer: drop the 14pt formatting on the FILENAME field runs
#  7. Normal style: overflowPunct false -> true

$d = $word.ActiveDocument

# --- 1. "Allo" -> "Al" -------------------------------------------------
$d.Content.Find.Execute("Allo", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Al", 2) | Out-Null

# --- 2. "Sportello unico per le attività produttive" -> "SUAP/SUE" -----
$d.Content.Find.Execute("Sportello unico per le attività produttive", `
                         $false, $true, $false, $false, $false, `
                         $true, 1, $false, "SUAP/SUE", 2) | Out-Null

# --- 5. Remove the "Si premette che il SUAP riceve ..." paragraph ------
#    (done before the generic " SUAP " substitution below so that this
#    paragraph's own "SUAP" text doesn't get rewritten first)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Si premette che il SUAP riceve*") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- 3. " SUAP " -> " SUAP/SUE " ---------------------------------------
$d.Content.Find.Execute(" SUAP ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " SUAP/SUE ", 2) | Out-Null

# --- 4. "ratica SUAP n°" -> "ratica SUAP/SUE n°" ------------------------
$d.Content.Find.Execute("ratica SUAP n" + [char]0xB0, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ratica SUAP/SUE n" + [char]0xB0, 2) | Out-Null

Write-Host "done"
